$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.722.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.445.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.99%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.77%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.444.64"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.75%  "

$ws.Range("E10").Value = "  +2.85%  "

$ws.Range("E11").Value = "  +2.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.25%  "

$ws.Range("E13").Value = "  +2.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000178"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.889.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.672.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.448.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.97%  "

$ws.Range("E19").Value = "  -1.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.60%  "

$ws.Range("E23").Value = "  +7.89%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "645.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.67%  "

$ws.Range("E27").Value = "  +17.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₆0552"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +98.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0987"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.05%  "

$ws.Range("E31").Value = "  +2.49%  "

$ws.Range("E32").Value = "  +2.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.88%  "

$ws.Range("E35").Value = "  +4.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.15%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("E38").Value = "  +3.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "153.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.55%  "

$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("E42").Value = "  +2.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.27%  "

$ws.Range("E45").Value = "  +1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "14.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +27.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "145.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.35%  "

$ws.Range("E51").Value = "  +2.56%  "
